$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 5: "2021年" in column A (mirrors the style of the existing year
# cells in column A, e.g. A4), followed by the data values for columns B..DK.
$row = 5

$ws.Cells.Item(4, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item($row, 1).Value = "2021年"

$values = @(
    -9.1, -4.8, $null, 27.1, -8.199999999999999, $null, -87.2, 31.8, 53.2, -76.8, $null, 21.1, 81.3,
    20, -98.90000000000001, $null, $null, $null, -83, -35.3, -52, -11.2, -8.800000000000001, 71, $null,
    3.2, -40.2, -59.3, -59.6, 120.4, 43.6, $null, $null, $null, $null, -14.6, 42.2, 47.1, -3.8, -19.9,
    -8.300000000000001, -77.59999999999999, $null, $null, $null, $null, -40.1, $null, 15.8, 84.3, 4.3,
    6, 121.3, -82.09999999999999, $null, -92.2, $null, -68, $null, $null, -34, $null, 5.9, $null,
    -11.7, -21.7, -70.3, $null, $null, -58.5, -29.8, $null, -32.9, -35.6, -52.4, -0.1, -99.7, -99.5,
    8.9, -95.8, $null, 137.3, -58, 206.7, 36.4, 100.1, $null, -80.40000000000001, -67.90000000000001,
    $null, 57.9, 12, 3.2, $null, $null, 19.6, -33.6, 86.90000000000001, -38.3, $null, 1.3, 15.9, 170.1,
    -40.3, $null, -57, -98.2, 164.8, -33.4, 198.4, 37.1, $null, 11, 42.9
)

$col = 2
foreach ($v in $values) {
    $ws.Cells.Item($row, $col).Value = $v
    $col = $col + 1
}
